# Fruta / hortaliza, semanal
# The underlying data rows (2-27) got re-shuffled: the contents of columns
# D (Fecha), I (Calidad), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), O (Origen) and P (Precio $/Kg) were
# redistributed among the rows. Everything else (Mercado, Region, Codreg,
# Categoria, Variedad, Unidad de comercializacion, Kg o Unidades,
# Clasificacion) stays put.
#
# row -> source row whose old D/I/J/K/L/M/O/P values now land on `row`
$mapping = @{
    2  = 22
    3  = 13
    4  = 26
    5  = 27
    6  = 25
    7  = 21
    8  = 10
    9  = 14
    10 = 7
    11 = 11
    12 = 5
    13 = 2
    14 = 8
    15 = 19
    16 = 24
    17 = 4
    18 = 20
    19 = 9
    20 = 12
    21 = 17
    22 = 3
    23 = 23
    24 = 15
    25 = 16
    26 = 18
    27 = 6
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "I", "J", "K", "L", "M", "O", "P")

# Snapshot the original values for every row before writing anything back,
# since several rows swap values with each other.
$snapshot = @{}
foreach ($row in $mapping.Keys) {
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Range("$col$row").Value()
    }
    $snapshot[$row] = $rowValues
}

foreach ($row in $mapping.Keys) {
    $srcRow = $mapping[$row]
    if ($srcRow -eq $row) {
        continue
    }
    $srcValues = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $srcValues[$col]
    }
}
